$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename shared string "ColCajicá" -> "ColCajica" (header for column CZ)
$ws.Range("CZ1").Value2 = "ColCajica"

# 2. Individual cell corrections (numeric <-> NaN placeholder text)
$ws.Range("AK33").Value2 = 1
$ws.Range("AK34").Value2 = 1
$ws.Range("AK35").Value2 = 1
$ws.Range("CI35").Value2 = "NaN"
$ws.Range("DM35").Value2 = "NaN"
$ws.Range("AK41").Value2 = 5
$ws.Range("AK42").Value2 = 5
$ws.Range("AK73").Value2 = 13
$ws.Range("AK74").Value2 = 13
$ws.Range("CM83").Value2 = "NaN"
$ws.Range("Q88").Value2 = "NaN"
$ws.Range("H90").Value2 = "NaN"
$ws.Range("AP92").Value2 = "NaN"
$ws.Range("CW95").Value2 = "NaN"
$ws.Range("CW96").Value2 = "NaN"
$ws.Range("AK100").Value2 = "NaN"
$ws.Range("CG130").Value2 = 100
$ws.Range("BY132").Value2 = "NaN"
$ws.Range("DQ140").Value2 = "NaN"
$ws.Range("DQ147").Value2 = "NaN"

# 3. Updated running-total values in column CN (rows 119-165)
$ws.Range("CN119").Value2 = 34
$ws.Range("CN120").Value2 = 34
$ws.Range("CN121").Value2 = 34
$ws.Range("CN122").Value2 = 34
$ws.Range("CN123").Value2 = 35
$ws.Range("CN124").Value2 = 36
$ws.Range("CN127").Value2 = 43
$ws.Range("CN128").Value2 = 45
$ws.Range("CN129").Value2 = 45
$ws.Range("CN130").Value2 = 54
$ws.Range("CN131").Value2 = 56
$ws.Range("CN132").Value2 = 63
$ws.Range("CN133").Value2 = 66
$ws.Range("CN134").Value2 = 73
$ws.Range("CN135").Value2 = 80
$ws.Range("CN137").Value2 = 104
$ws.Range("CN138").Value2 = 117
$ws.Range("CN139").Value2 = 167
$ws.Range("CN140").Value2 = 176
$ws.Range("CN141").Value2 = 183
$ws.Range("CN142").Value2 = 189
$ws.Range("CN143").Value2 = 200
$ws.Range("CN144").Value2 = 212
$ws.Range("CN145").Value2 = 230
$ws.Range("CN146").Value2 = 252
$ws.Range("CN147").Value2 = 269
$ws.Range("CN148").Value2 = 298
$ws.Range("CN149").Value2 = 335
$ws.Range("CN150").Value2 = 380
$ws.Range("CN151").Value2 = 403
$ws.Range("CN152").Value2 = 427
$ws.Range("CN153").Value2 = 448
$ws.Range("CN154").Value2 = 471
$ws.Range("CN155").Value2 = 476
$ws.Range("CN156").Value2 = 552
$ws.Range("CN157").Value2 = 566
$ws.Range("CN158").Value2 = 584
$ws.Range("CN159").Value2 = 603
$ws.Range("CN160").Value2 = 737
$ws.Range("CN161").Value2 = 790
$ws.Range("CN162").Value2 = 814
$ws.Range("CN163").Value2 = 841
$ws.Range("CN164").Value2 = 874
$ws.Range("CN165").Value2 = 892

# 4. Append new data row 167 (date 2020-08-18)
$row167 = @{
    "A" = 44061
    "B" = 489122
    "C" = 2669
    "D" = 64259
    "E" = 61571
    "F" = 171312
    "G" = 22244
    "H" = 2158
    "I" = 1737
    "J" = 4129
    "K" = 3435
    "L" = 6060
    "M" = 3524
    "N" = 15181
    "O" = 16087
    "P" = 3751
    "Q" = 2598
    "R" = 10549
    "S" = 4680
    "T" = 11763
    "U" = 7293
    "V" = 2181
    "W" = 657
    "X" = 3487
    "Y" = 10653
    "Z" = 9015
    "AA" = 4791
    "AB" = 39020
    "AC" = 691
    "AD" = 100
    "AE" = 154
    "AF" = 434
    "AG" = 20
    "AH" = 14
    "AI" = 198
    "AJ" = 1910
    "AK" = 1850
    "AL" = 34617
    "AM" = 5249
    "AN" = 2305
    "AO" = 30714
    "AP" = 765
    "AQ" = 18443
    "AR" = 1364
    "AS" = 5089
    "AT" = 1338
    "AU" = 1523
    "AV" = 2659
    "AW" = 1280
    "AX" = 925
    "AY" = 2438
    "AZ" = 2542
    "BA" = 36702
    "BB" = 9848
    "BC" = 1444
    "BD" = 6108
    "BE" = 2232
    "BF" = 272
    "BG" = 1346
    "BH" = 2447
    "BI" = 723
    "BJ" = 1857
    "BK" = 7154
    "BL" = 6508
    "BM" = 6173
    "BN" = 13358
    "BO" = 1829
    "BP" = 722
    "BQ" = 4348
    "BR" = 3925
    "BS" = 4357
    "BT" = 894
    "BU" = 1168
    "BV" = 1710
    "BW" = 1997
    "BX" = 483
    "BY" = 3531
    "BZ" = 1992
    "CA" = 791
    "CB" = 545
    "CC" = 1483
    "CD" = 1580
    "CE" = 706
    "CF" = 651
    "CG" = 3467
    "CH" = 862
    "CI" = 921
    "CJ" = 939
    "CK" = 1214
    "CL" = 1009
    "CM" = 957
    "CN" = 927
    "CO" = 887
    "CP" = 1000
    "CQ" = 438
    "CR" = 2686
    "CS" = 704
    "CT" = 705
    "CU" = 597
    "CV" = 975
    "CW" = 861
    "CX" = 514
    "CY" = 623
    "CZ" = 648
    "DA" = 953
    "DB" = 773
    "DC" = 858
    "DD" = 688
    "DE" = 309
    "DF" = 312
    "DG" = 565
    "DH" = 453
    "DI" = 356
    "DJ" = 500
    "DK" = 288
    "DL" = 488
    "DM" = 685
    "DN" = 494
    "DO" = 465
    "DP" = 339
    "DQ" = 506
    "DR" = 104175
    "DS" = 203347
    "DT" = 6154
    "DU" = 88445
    "DV" = 57742
    "DW" = 17995
    "DX" = 6330
}
foreach ($col in $row167.Keys) {
    $ws.Range("$col" + "167").Value2 = $row167[$col]
}

# 5. Update active selection to the newly added last cell
$ws.Range("DX167").Select()
